$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata": update Version/Date values and insert a new
# "Jurisdiction" property row (with an empty value) right after "Contact"
# and before "Description".
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# ---------------------------------------------------------------------------
# Sheet "Elements": content is unchanged; only the column "best fit" widths
# were recomputed (by the original publishing tool) and need updating to
# match. Hidden columns keep their hidden state.
# ---------------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# Columns 10, 12, 13, 14, 16-19 and 35 keep their original width (untouched
# below) - only the columns below were resized.
$els.Columns.Item(1).ColumnWidth = 14.0
$els.Columns.Item(2).ColumnWidth = 14.0
$els.Columns.Item(3).ColumnWidth = 9.0
$els.Columns.Item(4).ColumnWidth = 6.166666666666667
$els.Columns.Item(5).ColumnWidth = 4.5
$els.Columns.Item(6).ColumnWidth = 3.1666666666666665
$els.Columns.Item(7).ColumnWidth = 3.5
$els.Columns.Item(8).ColumnWidth = 11.833333333333334
$els.Columns.Item(9).ColumnWidth = 9.666666666666666
$els.Columns.Item(11).ColumnWidth = 7.5
$els.Columns.Item(15).ColumnWidth = 11.5
$els.Columns.Item(20).ColumnWidth = 7.0
$els.Columns.Item(21).ColumnWidth = 12.833333333333334
$els.Columns.Item(22).ColumnWidth = 13.166666666666666
$els.Columns.Item(23).ColumnWidth = 14.166666666666666
$els.Columns.Item(24).ColumnWidth = 13.833333333333334
$els.Columns.Item(25).ColumnWidth = 16.166666666666668
$els.Columns.Item(26).ColumnWidth = 14.333333333333334
$els.Columns.Item(27).ColumnWidth = 4.166666666666667
$els.Columns.Item(28).ColumnWidth = 17.166666666666668
$els.Columns.Item(29).ColumnWidth = 15.5
$els.Columns.Item(30).ColumnWidth = 12.666666666666666
$els.Columns.Item(31).ColumnWidth = 10.5
$els.Columns.Item(32).ColumnWidth = 14.166666666666666
$els.Columns.Item(33).ColumnWidth = 7.333333333333333
$els.Columns.Item(34).ColumnWidth = 7.666666666666667

$els.Columns.Item(3).Hidden = $true
$els.Columns.Item(4).Hidden = $true
$els.Columns.Item(31).Hidden = $true
$els.Columns.Item(32).Hidden = $true
$els.Columns.Item(33).Hidden = $true
